# Auto-generated Excel COM-interop script to apply updated TPM values
# to the Cx3cl1-Itgav LR-pairs worksheet (commit: "update scripts wuth new tpm")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 9.912502333333334
$ws.Range("H2").Value = 29.737507
$ws.Range("I2").Value = 0.306800202315277
$ws.Range("J2").Value = 0.3105483022825659
$ws.Range("M2").Value = 13.89934866666667
$ws.Range("N2").Value = 41.69804600000001
$ws.Range("O2").Value = 0.04853507553134179
$ws.Range("P2").Value = 0.04999273878390351
$ws.Range("Q2").Value = 137.7773260901469
$ws.Range("R2").Value = 1239.995934811322
$ws.Range("S2").Value = 0.01489057099240291
$ws.Range("T2").Value = 0.01552516015579702
# Row 3
$ws.Range("G3").Value = 9.912502333333334
$ws.Range("H3").Value = 29.737507
$ws.Range("I3").Value = 0.306800202315277
$ws.Range("J3").Value = 0.3105483022825659
$ws.Range("O3").Value = 0.245697991654417
$ws.Range("P3").Value = 0.253077086664408
$ws.Range("Q3").Value = 697.466975074654
$ws.Range("R3").Value = 6277.202775671885
$ws.Range("S3").Value = 0.0753801935480324
$ws.Range("T3").Value = 0.07859265961024968
# Row 4
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = 9.912502333333334
$ws.Range("H4").Value = 29.737507
$ws.Range("I4").Value = 0.306800202315277
$ws.Range("J4").Value = 0.3105483022825659
$ws.Range("M4").Value = 82.007665
$ws.Range("N4").Value = 246.022995
$ws.Range("O4").Value = 0.2863622109480123
$ws.Range("P4").Value = 0.2949625822722868
$ws.Range("Q4").Value = 812.9011706637184
$ws.Range("R4").Value = 7316.110535973466
$ws.Range("S4").Value = 0.08785598425430023
$ws.Range("T4").Value = 0.09160012916154031
# Row 5
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 9.912502333333334
$ws.Range("H5").Value = 29.737507
$ws.Range("I5").Value = 0.306800202315277
$ws.Range("J5").Value = 0.3105483022825659
$ws.Range("M5").Value = 25.0501465
$ws.Range("N5").Value = 50.100293
$ws.Range("O5").Value = 0.0874724982879541
$ws.Range("P5").Value = 0.06006638442832619
$ws.Range("Q5").Value = 248.3096356315918
$ws.Range("R5").Value = 1489.857813789551
$ws.Range("S5").Value = 0.02683658017176704
$ws.Range("T5").Value = 0.01865351370846865
# Row 6
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("G6").Value = 9.912502333333334
$ws.Range("H6").Value = 29.737507
$ws.Range("I6").Value = 0.306800202315277
$ws.Range("J6").Value = 0.3105483022825659
$ws.Range("M6").Value = 95.05788666666668
$ws.Range("N6").Value = 285.17366
$ws.Range("O6").Value = 0.3319322235782747
$ws.Range("P6").Value = 0.3419012078510756
$ws.Range("Q6").Value = 942.261523385069
$ws.Range("R6").Value = 8480.353710465621
$ws.Range("S6").Value = 0.1018368733487744
$ws.Range("T6").Value = 0.1061768396465102
# Row 7
$ws.Range("I7").Value = 0.648195099606228
$ws.Range("J7").Value = 0.6561139341222959
$ws.Range("M7").Value = 13.89934866666667
$ws.Range("N7").Value = 41.69804600000001
$ws.Range("O7").Value = 0.04853507553134179
$ws.Range("P7").Value = 0.04999273878390351
$ws.Range("Q7").Value = 291.0903804317196
$ws.Range("R7").Value = 2619.813423885476
$ws.Range("S7").Value = 0.03146019811843389
$ws.Range("T7").Value = 0.03280093252105521
# Row 8
$ws.Range("I8").Value = 0.648195099606228
$ws.Range("J8").Value = 0.6561139341222959
$ws.Range("O8").Value = 0.245697991654417
$ws.Range("P8").Value = 0.253077086664408
$ws.Range("S8").Value = 0.159260234173485
$ws.Range("T8").Value = 0.166047402967594
# Row 9
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("I9").Value = 0.648195099606228
$ws.Range("J9").Value = 0.6561139341222959
$ws.Range("M9").Value = 82.007665
$ws.Range("N9").Value = 246.022995
$ws.Range("O9").Value = 0.2863622109480123
$ws.Range("P9").Value = 0.2949625822722868
$ws.Range("Q9").Value = 1717.464823399663
$ws.Range("R9").Value = 15457.18341059697
$ws.Range("S9").Value = 0.1856185818489065
$ws.Range("T9").Value = 0.1935290602735414
# Row 10
$ws.Range("D10").Value = "MuSCs"
$ws.Range("I10").Value = 0.648195099606228
$ws.Range("J10").Value = 0.6561139341222959
$ws.Range("M10").Value = 25.0501465
$ws.Range("N10").Value = 50.100293
$ws.Range("O10").Value = 0.0874724982879541
$ws.Range("P10").Value = 0.06006638442832619
$ws.Range("Q10").Value = 524.6185882107263
$ws.Range("R10").Value = 3147.711529264358
$ws.Range("S10").Value = 0.05669924474056601
$ws.Range("T10").Value = 0.03941039179577131
# Row 11
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("I11").Value = 0.648195099606228
$ws.Range("J11").Value = 0.6561139341222959
$ws.Range("M11").Value = 95.05788666666668
$ws.Range("N11").Value = 285.17366
$ws.Range("O11").Value = 0.3319322235782747
$ws.Range("P11").Value = 0.3419012078510756
$ws.Range("Q11").Value = 1990.772161805996
$ws.Range("R11").Value = 17916.94945625396
$ws.Range("S11").Value = 0.2151568407248365
$ws.Range("T11").Value = 0.224326146564334
# Row 12
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.169852
$ws.Range("H12").Value = 2.339704
$ws.Range("I12").Value = 0.03620789364881174
$ws.Range("J12").Value = 0.0244334908452053
$ws.Range("M12").Value = 13.89934866666667
$ws.Range("N12").Value = 41.69804600000001
$ws.Range("O12").Value = 0.04853507553134179
$ws.Range("P12").Value = 0.04999273878390351
$ws.Range("Q12").Value = 16.26018083639734
$ws.Range("R12").Value = 97.56108501838402
$ws.Range("S12").Value = 0.001757352853075868
$ws.Range("T12").Value = 0.001221497125403246
# Row 13
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.169852
$ws.Range("H13").Value = 2.339704
$ws.Range("I13").Value = 0.03620789364881174
$ws.Range("J13").Value = 0.0244334908452053
$ws.Range("O13").Value = 0.245697991654417
$ws.Range("P13").Value = 0.253077086664408
$ws.Range("Q13").Value = 82.31353782195335
$ws.Range("R13").Value = 493.8812269317201
$ws.Range("S13").Value = 0.008896206751549765
$ws.Range("T13").Value = 0.00618355668014604
# Row 14
$ws.Range("D14").Value = "Inflammatory-Mac"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1.169852
$ws.Range("H14").Value = 2.339704
$ws.Range("I14").Value = 0.03620789364881174
$ws.Range("J14").Value = 0.0244334908452053
$ws.Range("M14").Value = 82.007665
$ws.Range("N14").Value = 246.022995
$ws.Range("O14").Value = 0.2863622109480123
$ws.Range("P14").Value = 0.2949625822722868
$ws.Range("Q14").Value = 95.93683091558002
$ws.Range("R14").Value = 575.62098549348
$ws.Range("S14").Value = 0.01036857247904422
$ws.Range("T14").Value = 0.007206965553628033
# Row 15
$ws.Range("D15").Value = "MuSCs"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1.169852
$ws.Range("H15").Value = 2.339704
$ws.Range("I15").Value = 0.03620789364881174
$ws.Range("J15").Value = 0.0244334908452053
$ws.Range("M15").Value = 25.0501465
$ws.Range("N15").Value = 50.100293
$ws.Range("O15").Value = 0.0874724982879541
$ws.Range("P15").Value = 0.06006638442832619
$ws.Range("Q15").Value = 29.304963983318
$ws.Range("R15").Value = 117.219855933272
$ws.Range("S15").Value = 0.003167194915206109
$ws.Range("T15").Value = 0.00146763145403409
# Row 16
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1.169852
$ws.Range("H16").Value = 2.339704
$ws.Range("I16").Value = 0.03620789364881174
$ws.Range("J16").Value = 0.0244334908452053
$ws.Range("M16").Value = 95.05788666666668
$ws.Range("N16").Value = 285.17366
$ws.Range("O16").Value = 0.3319322235782747
$ws.Range("P16").Value = 0.3419012078510756
$ws.Range("Q16").Value = 111.2036588327734
$ws.Range("R16").Value = 667.2219529966402
$ws.Range("S16").Value = 0.01201856664993577
$ws.Range("T16").Value = 0.00835384003199389
# Row 17
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.2842186666666667
$ws.Range("H17").Value = 0.852656
$ws.Range("I17").Value = 0.00879680442968319
$ws.Range("J17").Value = 0.008904272749933054
$ws.Range("M17").Value = 13.89934866666667
$ws.Range("N17").Value = 41.69804600000001
$ws.Range("O17").Value = 0.04853507553134179
$ws.Range("P17").Value = 0.04999273878390351
$ws.Range("Q17").Value = 3.950454345575112
$ws.Range("R17").Value = 35.554089110176
$ws.Range("S17").Value = 0.0004269535674291156
$ws.Range("T17").Value = 0.0004451489816480333
# Row 18
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("G18").Value = 0.2842186666666667
$ws.Range("H18").Value = 0.852656
$ws.Range("I18").Value = 0.00879680442968319
$ws.Range("J18").Value = 0.008904272749933054
$ws.Range("O18").Value = 0.245697991654417
$ws.Range("P18").Value = 0.253077086664408
$ws.Range("Q18").Value = 19.99829377423111
$ws.Range("R18").Value = 179.98464396808
$ws.Range("S18").Value = 0.002161357181349839
$ws.Range("T18").Value = 0.002253467406418334
# Row 19
$ws.Range("D19").Value = "Inflammatory-Mac"
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 0.2842186666666667
$ws.Range("H19").Value = 0.852656
$ws.Range("I19").Value = 0.00879680442968319
$ws.Range("J19").Value = 0.008904272749933054
$ws.Range("M19").Value = 82.007665
$ws.Range("N19").Value = 246.022995
$ws.Range("O19").Value = 0.2863622109480123
$ws.Range("P19").Value = 0.2949625822722868
$ws.Range("Q19").Value = 23.30810920274667
$ws.Range("R19").Value = 209.77298282472
$ws.Range("S19").Value = 0.002519072365761347
$ws.Range("T19").Value = 0.00262642728357701
# Row 20
$ws.Range("D20").Value = "MuSCs"
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 0.3333333333333333
$ws.Range("G20").Value = 0.2842186666666667
$ws.Range("H20").Value = 0.852656
$ws.Range("I20").Value = 0.00879680442968319
$ws.Range("J20").Value = 0.008904272749933054
$ws.Range("M20").Value = 25.0501465
$ws.Range("N20").Value = 50.100293
$ws.Range("O20").Value = 0.0874724982879541
$ws.Range("P20").Value = 0.06006638442832619
$ws.Range("Q20").Value = 7.119719238034667
$ws.Range("R20").Value = 42.718315428208
$ws.Range("S20").Value = 0.0007694784604149299
$ws.Range("T20").Value = 0.000534847470052148
# Row 21
$ws.Range("D21").Value = "Resolving-Mac"
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 0.3333333333333333
$ws.Range("G21").Value = 0.2842186666666667
$ws.Range("H21").Value = 0.852656
$ws.Range("I21").Value = 0.00879680442968319
$ws.Range("J21").Value = 0.008904272749933054
$ws.Range("M21").Value = 95.05788666666668
$ws.Range("N21").Value = 285.17366
$ws.Range("O21").Value = 0.3319322235782747
$ws.Range("P21").Value = 0.3419012078510756
$ws.Range("Q21").Value = 27.01722580455111
$ws.Range("R21").Value = 243.15503224096
$ws.Range("S21").Value = 0.002919942854727958
$ws.Range("T21").Value = 0.01865351370846865
